$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New activity-log rows appended below the existing A1:D353 data block
# (User, Activity, Date, Time) -- mirrors the 05/13/24-05/14/24 LOG-IN entries.
$newRows = @(
  @(354, "q", "LOG-IN", "05/13/24", "14:58:11"),
  @(355, "q", "LOG-IN", "05/13/24", "14:59:03"),
  @(356, "q", "LOG-IN", "05/13/24", "15:04:37"),
  @(357, "q", "LOG-IN", "05/13/24", "16:11:09"),
  @(358, "q", "LOG-IN", "05/13/24", "16:11:29"),
  @(359, "q", "LOG-IN", "05/13/24", "16:11:46"),
  @(360, "q", "LOG-IN", "05/13/24", "16:12:56"),
  @(361, "q", "LOG-IN", "05/13/24", "16:13:20"),
  @(362, "q", "LOG-IN", "05/13/24", "16:14:11"),
  @(363, "q", "LOG-IN", "05/14/24", "14:17:00"),
  @(364, "q", "LOG-IN", "05/14/24", "14:18:04"),
  @(365, "q", "LOG-IN", "05/14/24", "14:35:45"),
  @(366, "q", "LOG-IN", "05/14/24", "14:36:41"),
  @(367, "q", "LOG-IN", "05/14/24", "14:43:32"),
  @(368, "q", "LOG-IN", "05/14/24", "14:51:07"),
  @(369, "q", "LOG-IN", "05/14/24", "15:43:20"),
  @(370, "q", "LOG-IN", "05/14/24", "15:49:01"),
  @(371, "q", "LOG-IN", "05/14/24", "21:29:16"),
  @(372, "q", "LOG-IN", "05/14/24", "21:33:34"),
  @(373, "q", "LOG-IN", "05/14/24", "21:41:33"),
  @(374, "q", "LOG-IN", "05/14/24", "21:46:34"),
  @(375, "q", "LOG-IN", "05/14/24", "21:47:18"),
  @(376, "q", "LOG-IN", "05/14/24", "21:51:30"),
  @(377, "q", "LOG-IN", "05/14/24", "21:54:11"),
  @(378, "q", "LOG-IN", "05/14/24", "21:55:23"),
  @(379, "q", "LOG-IN", "05/14/24", "22:00:22"),
  @(380, "q", "LOG-IN", "05/14/24", "22:00:56"),
  @(381, "q", "LOG-IN", "05/14/24", "22:07:17"),
  @(382, "q", "LOG-IN", "05/14/24", "22:08:14"),
  @(383, "q", "LOG-IN", "05/14/24", "22:08:59"),
  @(384, "q", "LOG-IN", "05/14/24", "22:24:02"),
  @(385, "q", "LOG-IN", "05/14/24", "22:24:40"),
  @(386, "q", "LOG-IN", "05/14/24", "22:25:22")
)

$firstRow = $newRows[0][0]
$lastRow = $newRows[$newRows.Count - 1][0]

# Format as Text first so date- and time-looking strings ("05/13/24", "14:58:11")
# are stored verbatim instead of being auto-converted to date/time serials --
# matches the source data which is plain text in every column.
$targetRange = $ws.Range("A$($firstRow):D$($lastRow)")
$targetRange.NumberFormat = "@"

foreach ($row in $newRows) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
}

# Drop the temporary Text format again so the new cells keep the workbook-wide
# default (General) style, same as every other row in the sheet.
$targetRange.ClearFormats()

